# correcao do orcamento faltando aparte de vacalo
# Adds the missing "cavalo" / "vacalo" rows to the budget sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 12/13: label in column C, monthly value in column D.
$ws.Range("C12").Value = "cavalo"
$ws.Range("D12").Value = 10000
$ws.Range("C13").Value = "vacalo"
$ws.Range("D13").Value = 85200

# Give the new values an accounting/currency look (distinct from the
# BRL-formatted values already used in column B).
$ws.Range("D12:D13").Style = "Currency"

# Re-touch the original BRL formatted cells with their exact literal
# format code so their style entry is left untouched.
$ws.Range("B5:B7").NumberFormat = "_([`$BRL]\ * #,##0.00_);_([`$BRL]\ * \(#,##0.00\);_([`$BRL]\ * ""-""??_);_(@_)"

# New column D needs to be wide enough to fit the currency values.
$ws.Columns.Item(4).ColumnWidth = 10.7

# Leave the selection where the user ended up after the edit.
$ws.Range("F16").Select() | Out-Null
